$d = $word.ActiveDocument
$tbl = $d.Tables(1)
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellParaXml($cell, $inner) {
    $p = $cell.Range.Paragraphs(1)
    $xml = '<w:p ' + $wns + '><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' + $inner + '</w:p>'
    $p.Range.InsertXML($xml)
}

# Row 1, Col 1: "Short title of decisión" -> spell-checked run split
$inner1 = '<w:r><w:t xml:space="preserve">Short </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/><w:r><w:t>title</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/><w:r><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> decisión</w:t></w:r>'
Set-CellParaXml $tbl.Cell(1,1) $inner1

# Row 3, Col 2: remove the _GoBack bookmark that sits between "1" and "/10/2019"
$inner3 = '<w:r><w:t>3</w:t></w:r>' +
          '<w:r><w:t>1</w:t></w:r>' +
          '<w:r><w:t>/10/2019</w:t></w:r>'
Set-CellParaXml $tbl.Cell(3,2) $inner3

# Row 4, Col 1: "Creator of decisión" -> spell-checked run split
$inner4 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Creator</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/><w:r><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> decisión</w:t></w:r>'
Set-CellParaXml $tbl.Cell(4,1) $inner4

# Row 5, Col 1: "Description" -> wrapped with spell-check markers
$inner5 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Description</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-CellParaXml $tbl.Cell(5,1) $inner5

# Row 6, Col 1: "Decision's Rational" -> spell-checked run split
$inner6a = '<w:proofErr w:type="spellStart"/><w:r><w:t>Decision' + [char]8217 + 's</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>Rational</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-CellParaXml $tbl.Cell(6,1) $inner6a

# Row 6, Col 2: split " secure, " into " " + "secure" + ", ", preserve surrounding runs
$inner6b = '<w:r><w:t>Se propone el uso de radios militares universales DSP 9000</w:t></w:r>' +
           '<w:r><w:t>, de la empresa T</w:t></w:r>' +
           '<w:r><w:t>CC</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>secure</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
           '<w:r><w:t xml:space="preserve">para la comunicación con la policía y </w:t></w:r>' +
           '<w:r><w:t>la UME. Esta</w:t></w:r>' +
           '<w:r><w:t>s</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> pueden establecer comunicaciones por canales concretos y cifrados.</w:t></w:r>'
Set-CellParaXml $tbl.Cell(6,2) $inner6b

# Row 8, Col 1: "Requirements (decisión drivers)" -> spell-checked run split
$inner8 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Requirements</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> (decisión drivers)</w:t></w:r>'
Set-CellParaXml $tbl.Cell(8,1) $inner8

# Row 9, Col 1: "Alternative decisions (options)" -> spell-checked run split
$inner9 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Alternative</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/><w:r><w:t>decisions</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/><w:r><w:t>options</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>)</w:t></w:r>'
Set-CellParaXml $tbl.Cell(9,1) $inner9

# Row 10, Col 1: "Decision outcome (options selected)" -> spell-checked run split
$inner10a = '<w:proofErr w:type="spellStart"/><w:r><w:t>Decision</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>outcome</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>options</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>selected</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t>)</w:t></w:r>'
Set-CellParaXml $tbl.Cell(10,1) $inner10a

# Row 10, Col 2: empty -> add "?"
$inner10b = '<w:r><w:t>?</w:t></w:r>'
Set-CellParaXml $tbl.Cell(10,2) $inner10b

# Row 11, Col 2: empty -> add "?"
$inner11b = '<w:r><w:t>?</w:t></w:r>'
Set-CellParaXml $tbl.Cell(11,2) $inner11b

# Row 12, Col 1: "Cons opciones" -> spell-checked run split
$inner12a = '<w:proofErr w:type="spellStart"/><w:r><w:t>Cons</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> opciones</w:t></w:r>'
Set-CellParaXml $tbl.Cell(12,1) $inner12a

# Row 12, Col 2: empty w/ underline paragraph mark -> drop underline, add "?"
$p12b = $tbl.Cell(12,2).Range.Paragraphs(1)
$xml12b = '<w:p ' + $wns + '><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>?</w:t></w:r></w:p>'
$p12b.Range.InsertXML($xml12b)

# Row 13, Col 1: "Link to other decisions" -> spell-checked run split
$inner13 = '<w:r><w:t xml:space="preserve">Link </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>other</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>decisions</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-CellParaXml $tbl.Cell(13,1) $inner13

# Row 14, Col 1: "Link to architecture artifacts" -> spell-checked run split
$inner14 = '<w:r><w:t xml:space="preserve">Link </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>architecture</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>artifacts</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-CellParaXml $tbl.Cell(14,1) $inner14

# Final paragraph after the table: add the _GoBack bookmark there (moved from the Date cell)
$last = $d.Paragraphs($d.Paragraphs.Count)
$xmlLast = '<w:p ' + $wns + '><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$last.Range.InsertXML($xmlLast)
